# Insert a new data row at row 43 (pushing the existing rows 43..110 down to
# 44..111) and populate it with the new "Madrigal" / "Primera" record.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("43:43").Insert()

$ws.Cells.Item(43, 1).Value  = 5
$ws.Cells.Item(43, 2).Value  = "Macroferia Regional de Talca"
$ws.Cells.Item(43, 3).Value  = "Maule"
$ws.Cells.Item(43, 4).Value  = 44868
$ws.Cells.Item(43, 5).Value  = 7
$ws.Cells.Item(43, 6).Value  = 100112013
$ws.Cells.Item(43, 7).Value  = "Alcachofa"
$ws.Cells.Item(43, 8).Value  = "Madrigal"
$ws.Cells.Item(43, 9).Value  = "Primera"
$ws.Cells.Item(43, 10).Value = 200
$ws.Cells.Item(43, 11).Value = 8000
$ws.Cells.Item(43, 12).Value = 8000
$ws.Cells.Item(43, 13).Value = 8000
$ws.Cells.Item(43, 14).Value = "$/caja 40 unidades"
$ws.Cells.Item(43, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(43, 16).Value = 200
$ws.Cells.Item(43, 17).Value = 40
$ws.Cells.Item(43, 18).Value = "Hortaliza"
